# Commit: Added new results for evaluating the impact of code superpages on dtlb
#
# Adds a new first worksheet "data-superpg" with a small comparison table of
# perf-counter data across four superpage configurations, plus derived ratio
# columns.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new sheet as the first tab -----------------------------
$ws = $wb.Worksheets.Add()
$ws.Name = "data-superpg"
# Worksheets.Add() inserts before the (previously) active sheet, which is
# already the first tab, so $ws now sits at index 1.

$ws.Tab.ColorIndex = $ws.Tab.ColorIndex  # no-op touch, keeps default tab

# --- 2. Intro lines ---------------------------------------------------------
$ws.Range("A1").Value2 = "iterations = 100"
$ws.Range("A2").Value2 = "4 hyperthreads/2cores"
$ws.Range("A3").Value2 = " node index.js >& /dev/null"

# --- 3. Column header row (row 5) ------------------------------------------
$headers = @{
  "B5" = "code no superpage, data no superpage";
  "C5" = "code no superpage, data superpage";
  "D5" = "code superpage, data no superpage";
  "E5" = "code superpage, data superpage";
  "G5" = "Col C/Col B";
  "H5" = "Col D/ Col B";
  "I5" = "Col E/ Col B";
}
foreach ($addr in $headers.Keys) {
  $ws.Range($addr).Value2 = $headers[$addr]
}
$ws.Range("B5:E5,G5:I5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 78.75

# --- 4. Data rows 6-18 -------------------------------------------------------
# row, label, colB, colC, colD, colE
$rows = @(
  @(6,  "CPU_CLK_UNHALTED.THREAD_P",            35388438420.707397, 35430393799.9049,   34921808041.769897, 34880732666.752502),
  @(7,  "DTLB_LOAD_MISSES.WALK_COMPLETED",      1378851.02,         1383086.2875000001, 1184337.17,          1119978.6399999999),
  @(8,  "DTLB_LOAD_MISSES.WALK_PENDING",        188655874.54249999, 189362021.36750001, 176016859.72749999,  168905688.095),
  @(9,  "DTLB_STORE_MISSES.WALK_COMPLETED",     1373329.8425,       1370859.4975000001, 1364951.7075,        1341041.165),
  @(10, "DTLB_STORE_MISSES.WALK_PENDING",       83883996.337500006, 83678590.6875,       83413221.147499993, 80690930.282499999),
  @(11, "ITLB_MISSES.WALK_COMPLETED",           37412028.560000002, 37006602.630000003, 36050020.7575,       36679333.659999996),
  @(12, "ITLB_MISSES.WALK_PENDING",             826626855.07749999, 824567132.23249996, 779790397.3125,      791258657.28750002),
  @(13, "ICACHE_64B.IFTAG_STALL",                3748304983.8425002, 3716653125.1624999, 2711901250.4324999, 2723631579.2725),
  @(14, "CPU_CLK_UNHALTED.THREAD_P (os + usr)", 36030648327.915001, 36056316621.822502, 35520501426.997398, 35508158950.037399),
  @(15, "INST_RETIRED.ANY_P (os + user)",       46751655130.150002, 46762980871.297501, 46748528350.517502, 46747284103.184898),
  @(16, "INST_RETIRED.ANY_P",                   46086256393.125,    46097081223.3349,   46093250344.010002, 46091212580.652397),
  @(17, "CYCLE_ACTIVITY.STALLS_L3_MISS",        1177260374.2375,    1176961053.6724999, 1166643216.9300001, 1159155357.4875),
  @(18, "elapse time",                          1029.3789999999999, 1029.58,             1015.258,            1014.2805)
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Cells.Item($rowNum, 1).Value2 = $r[1]
  $ws.Cells.Item($rowNum, 2).Value2 = $r[2]
  $ws.Cells.Item($rowNum, 3).Value2 = $r[3]
  $ws.Cells.Item($rowNum, 4).Value2 = $r[4]
  $ws.Cells.Item($rowNum, 5).Value2 = $r[5]
  $ws.Cells.Item($rowNum, 7).Formula = "=C$rowNum/B$rowNum"
  $ws.Cells.Item($rowNum, 8).Formula = "=D$rowNum/B$rowNum"
  $ws.Cells.Item($rowNum, 9).Formula = "=E$rowNum/B$rowNum"
}

# Comma style (#,##0) on the data block
$ws.Range("B6:E18,G6:I18").Style = "Comma"

# Column A labels get a plain (no theme color) font; last row label is red
$ws.Range("A6:A18").Font.ThemeColor = 1
$ws.Range("A6:A17").Font.Name = "Calibri"
$ws.Range("A18").Font.Color = RGB(255,0,0)

# --- 5. Column widths / layout ---------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 36.8125
$ws.Range("B1:E1").ColumnWidth = 17.0625

$ws.Range("G22").Select()
Write-Output "done"
